$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Product Owner assignment (Amato / Amedeo)
# Clear the existing bold style first, then set values (new shared strings
# are appended, reusing existing ones where they already match).
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Value = "Amato"
$ws.Range("B4").Value = "Amedeo"

# Row 7: Scrum Master assignment rotates weekly - unbold the label cell and
# replace it with the new note text (adds a new shared string entry).
$ws.Range("A7").Font.Bold = $false
$ws.Range("A7").Value = "Jede Woche wechseln"

# Update the sheet's active selection to A7.
[void]$ws.Range("A7").Select()
